# Apply a cyan highlight to the runs (and paragraph mark) of the
# "(7 points) A world-oriented camera 'compass' ..." bullet group and the
# Audio-section bullet group, matching the authored diff.
$d = $word.ActiveDocument

$targetParagraphIndexes = @(56, 57, 58, 59, 64, 65, 66, 67, 68, 69)

foreach ($idx in $targetParagraphIndexes) {
    $p = $d.Paragraphs($idx)
    $p.Range.Font.HighlightColorIndex = 3
}

Write-Host "Applied cyan highlight to $($targetParagraphIndexes.Count) paragraphs"
